$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 6 "0 Results" rows at the Ad level (old rows 9,11,12,13,15,16), deleting bottom-up
$ws.Rows.Item(16).Delete()
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(9).Delete()

# Update header / label text for A1, A3, A7 (string content same, index renumbers automatically)
$ws.Range("A1").Value = "MARCH Campaign Level"
$ws.Range("A3").Value = "MARCH Ad Set Level"
$ws.Range("A7").Value = "MARCH Ad Level"

# Row 2 (Campaign-level) metric updates
$ws.Range("B2").Value = 583.39287734
$ws.Range("C2").Value = 47
$ws.Range("D2").Value = 63609
$ws.Range("E2").Value = 230
$ws.Range("F2").Value = 12.41261441148936
$ws.Range("G2").Value = 9.171546123032904
$ws.Range("H2").Value = 2.536490771043479
$ws.Range("I2").Value = 489.3617021276596
$ws.Range("J2").Value = 584.89944804
$ws.Range("K2").Value = 62
$ws.Range("L2").Value = 56540
$ws.Range("M2").Value = 245
$ws.Range("N2").Value = 9.433862065161291
$ws.Range("O2").Value = 10.34487881216838
$ws.Range("P2").Value = 2.387344685877551
$ws.Range("Q2").Value = 395.1612903225807
$ws.Range("R2").Value = -0.2575777264021228
$ws.Range("S2").Value = -24.19354838709678
$ws.Range("T2").Value = 31.57511023325678
$ws.Range("U2").Value = 12.50265298903431
$ws.Range("V2").Value = -11.34215983038213
$ws.Range("W2").Value = -6.122448979591836
$ws.Range("X2").Value = 6.247362856658605
$ws.Range("Y2").Value = 23.83847155883629
$ws.Range("Z2").Value = 106.7458255719991
$ws.Range("AA2").Value = -13.49332023745113
$ws.Range("AB2").Value = 53.61629419563344
$ws.Range("AC2").Value = 67.43044189852702

# Row 5 / Row 6 (Ad Set level) metric updates
$ws.Range("C5").Value = 416.25744999
$ws.Range("D5").Value = 34
$ws.Range("E5").Value = 43611
$ws.Range("F5").Value = 151
$ws.Range("G5").Value = 12.24286617617647
$ws.Range("H5").Value = 9.544781132971039
$ws.Range("I5").Value = 2.756671854238411
$ws.Range("J5").Value = 444.1176470588235
$ws.Range("C6").Value = 167.13542735
$ws.Range("D6").Value = 13
$ws.Range("E6").Value = 19998
$ws.Range("F6").Value = 79
$ws.Range("G6").Value = 12.85657133461538
$ws.Range("H6").Value = 8.35760712821282
$ws.Range("I6").Value = 2.115638320886076
$ws.Range("J6").Value = 607.6923076923076

# Rows 9-14 (Ad level) metric updates post-delete (B/C ad_name/ad_set_name already correct from shift)
$ws.Range("D9").Value = 98.19
$ws.Range("E9").Value = 9
$ws.Range("F9").Value = 11671
$ws.Range("G9").Value = 46
$ws.Range("H9").Value = 10.91
$ws.Range("I9").Value = 8.413160825978922
$ws.Range("J9").Value = 2.134565217391304
$ws.Range("K9").Value = 511.1111111111111
$ws.Range("D10").Value = 77.98773658
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 5245
$ws.Range("G10").Value = 28
$ws.Range("H10").Value = 9.7484670725
$ws.Range("I10").Value = 14.86896788941849
$ws.Range("J10").Value = 2.785276306428571
$ws.Range("K10").Value = 350
$ws.Range("D11").Value = 46.91542735
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 5487
$ws.Range("G11").Value = 15
$ws.Range("H11").Value = 23.457713675
$ws.Range("I11").Value = 8.550287470384546
$ws.Range("J11").Value = 3.127695156666667
$ws.Range("K11").Value = 750
$ws.Range("D12").Value = 338.26971341
$ws.Range("E12").Value = 26
$ws.Range("F12").Value = 38366
$ws.Range("G12").Value = 123
$ws.Range("H12").Value = 13.01037359269231
$ws.Range("I12").Value = 8.816913762445916
$ws.Range("J12").Value = 2.750160271626017
$ws.Range("K12").Value = 473.0769230769231
$ws.Range("D13").Value = 19.05
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 2687
$ws.Range("G13").Value = 17
$ws.Range("H13").Value = 19.05
$ws.Range("I13").Value = 7.089691105321921
$ws.Range("J13").Value = 1.120588235294118
$ws.Range("K13").Value = 1700
$ws.Range("D14").Value = 2.98
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 153
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 2.98
$ws.Range("I14").Value = 19.47712418300653
$ws.Range("J14").Value = 2.98
$ws.Range("K14").Value = 100
